$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 21 de Marzo de 2020 a las 11:16"

# Refresh the country data table (re-sorted by Casos totales desc, updated counts)
$ws.Cells.Item(7,1).Value = "Iran"
$ws.Cells.Item(7,2).Value = 20610
$ws.Cells.Item(7,3).Value = 966
$ws.Cells.Item(7,4).Value = 7635
$ws.Cells.Item(7,5).Value = 11419
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 123
$ws.Cells.Item(7,8).Value = 1556

$ws.Cells.Item(8,1).Value = "Alemania"
$ws.Cells.Item(8,2).Value = 20046
$ws.Cells.Item(8,3).Value = 198
$ws.Cells.Item(8,4).Value = 180
$ws.Cells.Item(8,5).Value = 19797
$ws.Cells.Item(8,6).Value = 2
$ws.Cells.Item(8,7).Value = 1
$ws.Cells.Item(8,8).Value = 69

$ws.Cells.Item(9,1).Value = "Estados Unidos"
$ws.Cells.Item(9,2).Value = 19774
$ws.Cells.Item(9,3).Value = 391
$ws.Cells.Item(9,4).Value = 147
$ws.Cells.Item(9,5).Value = 19352
$ws.Cells.Item(9,6).Value = 64
$ws.Cells.Item(9,7).Value = 19
$ws.Cells.Item(9,8).Value = 275

$ws.Cells.Item(15,1).Value = "Belgica"
$ws.Cells.Item(15,2).Value = 2815
$ws.Cells.Item(15,3).Value = 558
$ws.Cells.Item(15,4).Value = 204
$ws.Cells.Item(15,5).Value = 2544
$ws.Cells.Item(15,6).Value = 288
$ws.Cells.Item(15,7).Value = 30
$ws.Cells.Item(15,8).Value = 67

$ws.Cells.Item(16,1).Value = "Austria"
$ws.Cells.Item(16,2).Value = 2695
$ws.Cells.Item(16,3).Value = 46
$ws.Cells.Item(16,4).Value = 9
$ws.Cells.Item(16,5).Value = 2679
$ws.Cells.Item(16,6).Value = 14
$ws.Cells.Item(16,7).Value = 1
$ws.Cells.Item(16,8).Value = 7

$ws.Cells.Item(17,1).Value = "Noruega"
$ws.Cells.Item(17,2).Value = 1995
$ws.Cells.Item(17,3).Value = 36
$ws.Cells.Item(17,4).Value = 1
$ws.Cells.Item(17,5).Value = 1987
$ws.Cells.Item(17,6).Value = 27
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(17,8).Value = 7

$ws.Cells.Item(27,1).Value = "Israel"
$ws.Cells.Item(27,2).Value = 883
$ws.Cells.Item(27,3).Value = 178
$ws.Cells.Item(27,4).Value = 36
$ws.Cells.Item(27,5).Value = 846
$ws.Cells.Item(27,6).Value = 15
$ws.Cells.Item(27,7).Value = 0
$ws.Cells.Item(27,8).Value = 1

$ws.Cells.Item(28,1).Value = "Crucero"
$ws.Cells.Item(28,2).Value = 712
$ws.Cells.Item(28,3).Value = 0
$ws.Cells.Item(28,4).Value = 567
$ws.Cells.Item(28,5).Value = 137
$ws.Cells.Item(28,6).Value = 15
$ws.Cells.Item(28,7).Value = 0
$ws.Cells.Item(28,8).Value = 8

$ws.Cells.Item(31,1).Value = "Pakistan"
$ws.Cells.Item(31,2).Value = 534
$ws.Cells.Item(31,3).Value = 33
$ws.Cells.Item(31,4).Value = 13
$ws.Cells.Item(31,5).Value = 518
$ws.Cells.Item(31,6).Value = 0
$ws.Cells.Item(31,7).Value = 0
$ws.Cells.Item(31,8).Value = 3

$ws.Cells.Item(47,1).Value = "Estonia"
$ws.Cells.Item(47,2).Value = 306
$ws.Cells.Item(47,3).Value = 23
$ws.Cells.Item(47,4).Value = 1
$ws.Cells.Item(47,5).Value = 305
$ws.Cells.Item(47,6).Value = 1
$ws.Cells.Item(47,7).Value = 0
$ws.Cells.Item(47,8).Value = 0

$ws.Cells.Item(48,1).Value = "Barein"
$ws.Cells.Item(48,2).Value = 298
$ws.Cells.Item(48,3).Value = 0
$ws.Cells.Item(48,4).Value = 125
$ws.Cells.Item(48,5).Value = 172
$ws.Cells.Item(48,6).Value = 4
$ws.Cells.Item(48,7).Value = 0
$ws.Cells.Item(48,8).Value = 1

$ws.Cells.Item(49,1).Value = "Egipto"
$ws.Cells.Item(49,2).Value = 285
$ws.Cells.Item(49,3).Value = 0
$ws.Cells.Item(49,4).Value = 42
$ws.Cells.Item(49,5).Value = 235
$ws.Cells.Item(49,6).Value = 0
$ws.Cells.Item(49,7).Value = 0
$ws.Cells.Item(49,8).Value = 8

$ws.Cells.Item(81,1).Value = "Republica de Macedonia"
$ws.Cells.Item(81,2).Value = 77
$ws.Cells.Item(81,3).Value = 1
$ws.Cells.Item(81,4).Value = 1
$ws.Cells.Item(81,5).Value = 76
$ws.Cells.Item(81,6).Value = 1
$ws.Cells.Item(81,7).Value = 0
$ws.Cells.Item(81,8).Value = 0

$ws.Cells.Item(82,1).Value = "Albania"
$ws.Cells.Item(82,2).Value = 76
$ws.Cells.Item(82,3).Value = 6
$ws.Cells.Item(82,4).Value = 2
$ws.Cells.Item(82,5).Value = 72
$ws.Cells.Item(82,6).Value = 2
$ws.Cells.Item(82,7).Value = 0
$ws.Cells.Item(82,8).Value = 2

$ws.Cells.Item(83,1).Value = "Republica de Chipre"
$ws.Cells.Item(83,2).Value = 75
$ws.Cells.Item(83,3).Value = 0
$ws.Cells.Item(83,4).Value = 0
$ws.Cells.Item(83,5).Value = 75
$ws.Cells.Item(83,6).Value = 1
$ws.Cells.Item(83,7).Value = 0
$ws.Cells.Item(83,8).Value = 0

$ws.Cells.Item(84,1).Value = "Principado de Andorra"
$ws.Cells.Item(84,2).Value = 75
$ws.Cells.Item(84,3).Value = 0
$ws.Cells.Item(84,4).Value = 1
$ws.Cells.Item(84,5).Value = 74
$ws.Cells.Item(84,6).Value = 2
$ws.Cells.Item(84,7).Value = 0
$ws.Cells.Item(84,8).Value = 0

$ws.Cells.Item(85,1).Value = "Sri Lanka"
$ws.Cells.Item(85,2).Value = 73
$ws.Cells.Item(85,3).Value = 0
$ws.Cells.Item(85,4).Value = 3
$ws.Cells.Item(85,5).Value = 70
$ws.Cells.Item(85,6).Value = 0
$ws.Cells.Item(85,7).Value = 0
$ws.Cells.Item(85,8).Value = 0

$ws.Cells.Item(86,1).Value = "Republica Dominicana"
$ws.Cells.Item(86,2).Value = 72
$ws.Cells.Item(86,3).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 70
$ws.Cells.Item(86,6).Value = 0
$ws.Cells.Item(86,7).Value = 0
$ws.Cells.Item(86,8).Value = 2

$ws.Cells.Item(115,1).Value = "Ghana"
$ws.Cells.Item(115,2).Value = 19
$ws.Cells.Item(115,3).Value = 3
$ws.Cells.Item(115,4).Value = 0
$ws.Cells.Item(115,5).Value = 19
$ws.Cells.Item(115,6).Value = 0
$ws.Cells.Item(115,7).Value = 0
$ws.Cells.Item(115,8).Value = 0

$ws.Cells.Item(116,1).Value = "Jamaica"
$ws.Cells.Item(116,2).Value = 19
$ws.Cells.Item(116,3).Value = 0
$ws.Cells.Item(116,4).Value = 2
$ws.Cells.Item(116,5).Value = 16
$ws.Cells.Item(116,6).Value = 0
$ws.Cells.Item(116,7).Value = 0
$ws.Cells.Item(116,8).Value = 1

$ws.Cells.Item(117,1).Value = "Paraguay"
$ws.Cells.Item(117,2).Value = 18
$ws.Cells.Item(117,3).Value = 0
$ws.Cells.Item(117,4).Value = 0
$ws.Cells.Item(117,5).Value = 18
$ws.Cells.Item(117,6).Value = 1
$ws.Cells.Item(117,7).Value = 0
$ws.Cells.Item(117,8).Value = 0

$ws.Cells.Item(118,1).Value = "Ruanda"
$ws.Cells.Item(118,2).Value = 17
$ws.Cells.Item(118,3).Value = 0
$ws.Cells.Item(118,4).Value = 0
$ws.Cells.Item(118,5).Value = 17
$ws.Cells.Item(118,6).Value = 0
$ws.Cells.Item(118,7).Value = 0
$ws.Cells.Item(118,8).Value = 0

$ws.Cells.Item(119,1).Value = "Macao"
$ws.Cells.Item(119,2).Value = 17
$ws.Cells.Item(119,3).Value = 0
$ws.Cells.Item(119,4).Value = 10
$ws.Cells.Item(119,5).Value = 7
$ws.Cells.Item(119,6).Value = 0
$ws.Cells.Item(119,7).Value = 0
$ws.Cells.Item(119,8).Value = 0

$ws.Cells.Item(124,1).Value = "Montenegro"
$ws.Cells.Item(124,2).Value = 14
$ws.Cells.Item(124,3).Value = 0
$ws.Cells.Item(124,4).Value = 0
$ws.Cells.Item(124,5).Value = 14
$ws.Cells.Item(124,6).Value = 0
$ws.Cells.Item(124,7).Value = 0
$ws.Cells.Item(124,8).Value = 0

$ws.Cells.Item(125,1).Value = "Puerto Rico"
$ws.Cells.Item(125,2).Value = 14
$ws.Cells.Item(125,3).Value = 0
$ws.Cells.Item(125,4).Value = 0
$ws.Cells.Item(125,5).Value = 14
$ws.Cells.Item(125,6).Value = 0
$ws.Cells.Item(125,7).Value = 0
$ws.Cells.Item(125,8).Value = 0

